$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 477, shifting existing rows 477-510 down to 478-511.
$ws.Rows.Item(477).Insert()

# Populate the newly inserted row 477 with the new data record.
$ws.Cells.Item(477, 1).Value = 3
$ws.Cells.Item(477, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(477, 3).Value = "Coquimbo"
$ws.Cells.Item(477, 4).Value2 = 44610
$ws.Cells.Item(477, 5).Value = 5
$ws.Cells.Item(477, 6).Value = 100112006
$ws.Cells.Item(477, 7).Value = "Repollo"
$ws.Cells.Item(477, 8).Value = "Crespo record"
$ws.Cells.Item(477, 9).Value = "Primera"
$ws.Cells.Item(477, 10).Value = 1200
$ws.Cells.Item(477, 11).Value = 1200
$ws.Cells.Item(477, 12).Value = 1300
$ws.Cells.Item(477, 13).Value = 1246
$ws.Cells.Item(477, 14).Value = "$/unidad"
$ws.Cells.Item(477, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(477, 16).Value = 1246
$ws.Cells.Item(477, 17).Value = 1
$ws.Cells.Item(477, 18).Value = "Hortaliza"
